$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 347.53333
$ws.Range("J33").Value = 474
$ws.Range("L33").Value = 474
$ws.Range("N33").Value = -932
$ws.Range("H43").Value = 4193.316
$ws.Range("J43").Value = 3854.5625
$ws.Range("L43").Value = 3854.5625
$ws.Range("N43").Value = -3992.5625
$ws.Range("H107").Value = 1493
$ws.Range("I107").Value = 1493
$ws.Range("K107").Value = 1493
$ws.Range("M107").Value = 427
$ws.Range("H132").Value = 26765.428
$ws.Range("I132").Value = 1840.4117
$ws.Range("J132").Value = 132696.75
$ws.Range("K132").Value = 5521.2351
$ws.Range("L132").Value = 398090.25
$ws.Range("M132").Value = -2991.2351
$ws.Range("N132").Value = -403150.25
$ws.Range("H137").Value = 5041
$ws.Range("I137").Value = 4317
$ws.Range("K137").Value = 12951
$ws.Range("M137").Value = -10401
$ws.Range("H138").Value = 3222.6216
$ws.Range("I138").Value = 1696
$ws.Range("J138").Value = 3578.8333
$ws.Range("K138").Value = 5088
$ws.Range("L138").Value = 10736.4999
$ws.Range("M138").Value = 52
$ws.Range("N138").Value = -21016.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1453.4706
$ws.Range("J88").Value = 1345.4286
$ws.Range("L88").Value = 1345.4286
$ws.Range("N88").Value = -2157.4286
$ws.Range("H91").Value = 1453.4706
$ws.Range("J91").Value = 1345.4286
$ws.Range("L91").Value = 1345.4286
$ws.Range("N91").Value = -4153.4286
$ws.Range("H102").Value = 14638.728
$ws.Range("I102").Value = 6780.778
$ws.Range("K102").Value = 6780.778
$ws.Range("M102").Value = -5158.778
$ws.Range("H130").Value = 170143
$ws.Range("J130").Value = 170143
$ws.Range("L130").Value = 170143
$ws.Range("N130").Value = -180183

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8518.071
$ws.Range("I20").Value = 5425.6
$ws.Range("K20").Value = 5425.6
$ws.Range("M20").Value = -5178.6
$ws.Range("H86").Value = 20041668
$ws.Range("I86").Value = 35787200
$ws.Range("J86").Value = 1899.3636
$ws.Range("K86").Value = 35787200
$ws.Range("L86").Value = 1899.3636
$ws.Range("M86").Value = -35786077
$ws.Range("N86").Value = -4145.3636
$ws.Range("H89").Value = 20041668
$ws.Range("I89").Value = 35787200
$ws.Range("J89").Value = 1899.3636
$ws.Range("K89").Value = 178936000
$ws.Range("L89").Value = 9496.817999999999
$ws.Range("M89").Value = -178930384
$ws.Range("N89").Value = -20728.818
$ws.Range("H105").Value = 2703
$ws.Range("I105").Value = 1490.4286
$ws.Range("J105").Value = 4825
$ws.Range("K105").Value = 1490.4286
$ws.Range("L105").Value = 4825
$ws.Range("M105").Value = 256.5714
$ws.Range("N105").Value = -8319
$ws.Range("H107").Value = 9193.521000000001
$ws.Range("I107").Value = 8756.277
$ws.Range("K107").Value = 8756.277
$ws.Range("M107").Value = -6836.277

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1797
$ws.Range("I31").Value = 1797
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1797
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1502
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1797
$ws.Range("I34").Value = 1797
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1797
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1595
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 2155.7
$ws.Range("I58").Value = 1365
$ws.Range("K58").Value = 1365
$ws.Range("M58").Value = -1162
$ws.Range("H105").Value = 2302.8333
$ws.Range("I105").Value = 2523.4
$ws.Range("K105").Value = 2523.4
$ws.Range("M105").Value = -776.4000000000001
$ws.Range("H132").Value = 2780
$ws.Range("I132").Value = 2112.2856
$ws.Range("J132").Value = 4338
$ws.Range("K132").Value = 6336.8568
$ws.Range("L132").Value = 13014
$ws.Range("M132").Value = -3806.8568
$ws.Range("N132").Value = -18074
$ws.Range("H136").Value = 2155.7
$ws.Range("I136").Value = 1365
$ws.Range("K136").Value = 4095
$ws.Range("M136").Value = -1545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 240.4375
$ws.Range("I10").Value = 52.666668
$ws.Range("J10").Value = 803.75
$ws.Range("K10").Value = 158.000004
$ws.Range("L10").Value = 2411.25
$ws.Range("M10").Value = -19.00000399999999
$ws.Range("N10").Value = -2689.25
$ws.Range("H34").Value = 1293.3334
$ws.Range("I34").Value = 148
$ws.Range("J34").Value = 2725
$ws.Range("K34").Value = 444
$ws.Range("L34").Value = 8175
$ws.Range("M34").Value = -360
$ws.Range("N34").Value = -8343
$ws.Range("H44").Value = 172
$ws.Range("I44").Value = 172
$ws.Range("K44").Value = 516
$ws.Range("M44").Value = -118
$ws.Range("H109").Value = 417.25
$ws.Range("I109").Value = 417.25
$ws.Range("K109").Value = 1251.75
$ws.Range("M109").Value = -211.75
$ws.Range("H137").Value = 1738.5
$ws.Range("I137").Value = 1707.4166
$ws.Range("J137").Value = 1925
$ws.Range("K137").Value = 5122.2498
$ws.Range("L137").Value = 5775
$ws.Range("M137").Value = -22.2497999999996
$ws.Range("N137").Value = -15975

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 27287
$ws.Range("J47").Value = 25930.5
$ws.Range("L47").Value = 25930.5
$ws.Range("N47").Value = -27066.5
$ws.Range("H48").Value = 31666.666
$ws.Range("J48").Value = 31666.666
$ws.Range("L48").Value = 31666.666
$ws.Range("N48").Value = -32636.666
$ws.Range("H99").Value = 38181.5
$ws.Range("I99").Value = 31363
$ws.Range("J99").Value = 45000
$ws.Range("K99").Value = 31363
$ws.Range("L99").Value = 45000
$ws.Range("M99").Value = -29117
$ws.Range("N99").Value = -49492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2677.6428
$ws.Range("H100").Value = 140388.5
$ws.Range("I100").Value = 140388.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 140388.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -139847.5
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4997.5
$ws.Range("J62").Value = 4995
$ws.Range("L62").Value = 4995
$ws.Range("N62").Value = -6243
$ws.Range("H65").Value = 4997.5
$ws.Range("J65").Value = 4995
$ws.Range("L65").Value = 24975
$ws.Range("N65").Value = -31215
$ws.Range("H135").Value = 44768.69
$ws.Range("J135").Value = 44768.69
$ws.Range("L135").Value = 44768.69
$ws.Range("N135").Value = -54908.69

